$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 8780
$ws.Range("F7").Value = 11196
$ws.Range("F13").Value = 125
$ws.Range("F15").Value = 304
$ws.Range("F16").Value = 261
$ws.Range("F18").Value = 88
$ws.Range("F20").Value = 422
$ws.Range("F22").Value = 1908
$ws.Range("F23").Value = 714
$ws.Range("F24").Value = 639
$ws.Range("F25").Value = 361
$ws.Range("F26").Value = 294
$ws.Range("F28").Value = 607
$ws.Range("F30").Value = 1312
$ws.Range("F33").Value = 7
$ws.Range("F36").Value = 466
$ws.Range("F38").Value = 357
$ws.Range("F40").Value = 34
$ws.Range("F44").Value = 120
$ws.Range("F45").Value = 816
$ws.Range("F46").Value = 659
$ws.Range("F48").Value = 164
$ws.Range("F49").Value = 148

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 58
$ws.Range("F18").Value = 68
$ws.Range("F19").Value = 109
$ws.Range("F23").Value = 53
$ws.Range("F24").Value = 59
$ws.Range("F25").Value = 396

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2840

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 8780
$ws.Range("F9").Value = 11196
$ws.Range("F14").Value = 125
$ws.Range("F15").Value = 304
$ws.Range("F18").Value = 422
$ws.Range("F19").Value = 1908
$ws.Range("F20").Value = 714
$ws.Range("F21").Value = 639
$ws.Range("F22").Value = 361
$ws.Range("F23").Value = 294
$ws.Range("F26").Value = 607
$ws.Range("F29").Value = 1312
$ws.Range("F37").Value = 357
$ws.Range("F42").Value = 120
$ws.Range("F44").Value = 53
$ws.Range("F45").Value = 396
$ws.Range("F46").Value = 659
$ws.Range("F48").Value = 164
$ws.Range("F49").Value = 148
